# Swap the contents of columns D ("codeforiati:group-code") and E
# ("codeforiati:group-name") for every used row on the sheet, including
# the header row. This matches the upstream change where the shared
# string table entries for each code/name pair were swapped, causing the
# two columns to exchange places in the rendered grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
